$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 714, shifting existing rows 714-755 down to 715-756.
$ws.Rows.Item(714).Insert()

# Populate the newly inserted row 714 with the new data point.
# Use a leading apostrophe so the date-like text is stored as plain text
# (matching the original inlineStr/text cells) rather than being parsed
# into a date serial value, then reset the style so no extra formatting
# (e.g. quote-prefix) is left behind on the cell.
$ws.Range("A714").Value = "'2026/01/29"
$ws.Range("A714").Style = "Normal"
$ws.Range("B714").Value = "木"
$ws.Range("C714").Value = 3
$ws.Range("D714").Value = 29
